$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("TC1")
$ws2 = $wb.Worksheets.Item("DragnDrop")

# --- Sheet1 (TC1): drop the 3 rows that only held now-unused locators
#     ("enter"/name/q, "typenumeric", click/xpath/Realeza) so the
#     remaining "click" + xpath(flow-news) row shifts up and becomes
#     row 5 (keeping its own original formatting). ---
$ws1.Rows("5:7").Delete()
$ws1.Range("A5").Value = 4

# --- Sheet2 (DragnDrop): the extra numeric value in J3 is no longer used ---
$ws2.Range("J3").ClearContents()

# --- Active sheet / selection bookkeeping: TC1 becomes the active tab ---
$ws2.Range("F3").Select()
$ws1.Select()
$ws1.Range("B7").Select()

Write-Output "done"
